$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.6712930202484151, -3.392556905746461, 2.365111112594603, 0.2393064647912979, -0.8791878223419189, -0.1872301995754242),
    @(0.5537151455879301, -4.361428594589236, 3.347476267814645, -0.1922698318958282, -0.9285151958465576, 0.8594874143600464),
    @(2.452674245834349, -4.409869003295896, 5.084140586853025, -3.570354700088501, -0.7802276611328125, -4.989242076873779),
    @(5.681596696376824, -3.525700151920317, 2.85166837722062, -1.221577763557434, 2.375196695327759, -2.503631114959717),
    @(1.737989616393984, -2.796510410308835, 1.393881118297584, 2.165364503860474, 0.5566509366035461, -0.4453207552433014),
    @(-5.58917605876923, -7.849099040031435, 6.400659620761871, 0.2727513313293457, 0.5925393104553223, 0.4948008358478546),
    @(2.166972637176496, -11.98566874265667, 9.183138275146462, -0.2121229618787765, 1.80510675907135, 1.96942949295044),
    @(-0.1177038192748863, -2.770210593938835, 3.860614097118379, 0.2755002379417419, 1.588096976280212, 2.037540912628174),
    @(3.127950906753536, -4.12096252441406, 3.36216964721679, 0.8869763612747192, 0.8231409192085266, 1.362993121147156),
    @(2.318384975194932, -3.449181020259855, 1.645497059822083, 0.2379320114850998, -0.7533495426177979, 0.1786780804395675),
    @(3.076034724712372, -2.935223340988159, 1.488467574119568, 0.6151412725448608, 1.230893492698669, -0.3686571717262268),
    @(4.307219874858854, -2.422872281074523, 1.113696080446245, 0.1403462886810302, 0.7915286421775818, 0.00137444678694),
    @(3.889325714111331, -2.331348705291748, 1.212168788909908, -0.3019201457500458, 0.041233405470848, -0.0345138870179653),
    @(4.565532350540161, -2.527925431728365, 0.07184810638427547, -0.1716531366109848, -0.2657263875007629, -0.0448985956609249),
    @(4.525413262844082, -2.978561615943909, -0.2339930295944185, -0.2078468948602676, 0.078801617026329, -0.026419922709465),
    @(3.761336714029313, -3.091410279273987, 0.3938730955123901, 0.0473420582711696, -0.0163406450301408, 0.0430659987032413),
    @(4.025719761848449, -2.884652662277222, 0.07831997871398896, -0.0447458773851394, 0.1368338167667389, -0.0710130855441093),
    @(3.922824805974961, -2.91312141418457, 0.03645527362823536, 0.0578794814646244, 0.06536258012056349, -0.0146607663482427),
    @(4.07339057922363, -2.916565942764282, 0.1658158779144287, -0.0494800843298435, 0.0968221426010131, 0.09865473955869671),
    @(3.406079268455506, -2.918866038322449, 0.1608816981315587, -0.08109235763549801, 0.040775254368782, 0.0303905457258224)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item(2 + $i, 3 + $j).Value = $row[$j]
    }
}
